$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.677.26"
$ws.Range("E2").Value = "  -2.41%  "
$ws.Range("D3").Value = "3.137.44"
$ws.Range("E3").Value = "  -8.25%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'564.56"
$ws.Range("E5").Value = "  -3.29%  "
$ws.Range("D6").Value = "'170.17"
$ws.Range("E6").Value = "  -4.74%  "
$ws.Range("D7").Value = "'0.616"
$ws.Range("E7").Value = "  -1.22%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").Value = "3.139.66"
$ws.Range("E9").Value = "  -8.05%  "
$ws.Range("E10").Value = "  -6.11%  "
$ws.Range("D11").Value = "'6.54"
$ws.Range("E11").Value = "  -6.03%  "
$ws.Range("D12").Value = "'0.393"
$ws.Range("E12").Value = "  -4.96%  "
$ws.Range("D13").Value = "3.680.58"
$ws.Range("E13").Value = "  -8.22%  "
$ws.Range("D14").Value = "'0.136"
$ws.Range("E14").Value = "  +0.85%  "
$ws.Range("D15").Value = "'27.00"
$ws.Range("D16").Value = "64.595.45"
$ws.Range("E16").Value = "  -2.59%  "
$ws.Range("E17").Value = "  -5.90%  "
$ws.Range("D18").Value = "3.136.37"
$ws.Range("E18").Value = "  -7.74%  "
$ws.Range("E19").Value = "  -4.10%  "
$ws.Range("E20").Value = "  -6.85%  "
$ws.Range("D21").Value = "'354.76"
$ws.Range("E21").Value = "  -3.41%  "
$ws.Range("E22").Value = "  -4.83%  "
$ws.Range("E23").Value = "  +0.39%  "
$ws.Range("D24").Value = "'68.36"
$ws.Range("E24").Value = "  -6.30%  "
$ws.Range("D25").Value = "'0.0000118"
$ws.Range("E25").Value = "  -6.99%  "
$ws.Range("E26").Value = "  -7.09%  "
$ws.Range("D27").Value = "'9.54"
$ws.Range("E27").Value = "  -2.95%  "
$ws.Range("E28").Value = "  -2.54%  "
$ws.Range("E29").Value = "  -0.10%  "
$ws.Range("E30").Value = "  -0.11%  "
$ws.Range("E31").Value = "  -5.05%  "
$ws.Range("E32").Value = "  -7.43%  "
$ws.Range("E33").Value = "  -6.49%  "
$ws.Range("D34").Value = "'6.63"
$ws.Range("E34").Value = "  -5.50%  "
$ws.Range("E35").Value = "  -5.45%  "
$ws.Range("D36").Value = "'1.44"
$ws.Range("E36").Value = "  -7.33%  "
$ws.Range("D37").Value = "'153.59"
$ws.Range("E37").Value = "  -5.70%  "
$ws.Range("D38").Value = "'0.830"
$ws.Range("E38").Value = "  -4.91%  "
$ws.Range("D39").Value = "'26.01"
$ws.Range("E39").Value = "  -5.75%  "
$ws.Range("E40").Value = "  -3.51%  "
$ws.Range("D41").Value = "'2.53"
$ws.Range("E41").Value = "  -2.33%  "
$ws.Range("D42").Value = "2.648.26"
$ws.Range("E42").Value = "  -1.97%  "
$ws.Range("E43").Value = "  -6.56%  "
$ws.Range("D44").Value = "'6.00"
$ws.Range("E44").Value = "  -5.00%  "
$ws.Range("D45").Value = "'24.16"
$ws.Range("E45").Value = "  -4.19%  "
$ws.Range("B46").Value = "Hedera"
$ws.Range("C46").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D46").Value = "'0.0653"
$ws.Range("E46").Value = "  -5.18%  "
$ws.Range("B47").Value = "OKB"
$ws.Range("C47").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D47").Value = "'39.00"
$ws.Range("E47").Value = "  -2.30%  "
$ws.Range("D48").Value = "'320.40"
$ws.Range("E48").Value = "  -4.87%  "
$ws.Range("E49").Value = "  -4.89%  "
$ws.Range("E50").Value = "  -2.49%  "
$ws.Range("D51").Value = "'1.00"
$ws.Range("E51").Value = "  +0.03%  "
